$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row above row 3 for the new "Testing" scenario.
#    (Excel's default Insert behaviour copies the format of the row above,
#    which we then overwrite with the correct per-cell formatting below.)
$ws.Rows("3").Insert()

# 2. Fill in the values for the new row 3 ("Testing" scenario).
$ws.Range("A3").Value = "Testing"
$ws.Range("C3").Value = "YES"
$ws.Range("D3").Value = "productCatalogPage"
$ws.Range("E3").Value = "productDetailPage"
$ws.Range("F3").Value = "productCatalogPage"
$ws.Range("G3").Value = "productDetailPage"
$ws.Range("H3").Value = "productCatalogPage"
$ws.Range("I3").Value = "productDetailPage"
$ws.Range("J3").Value = "cartCheck"

# 3. Match formatting for row 3 from existing, already-correctly-styled cells.
$ws.Range("A2:C2").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)
$ws.Range("D5:E5").Copy()
$ws.Range("D3:E3").PasteSpecial(-4122)
$ws.Range("D3:E3").Copy()
$ws.Range("F3:G3").PasteSpecial(-4122)
$ws.Range("H3:I3").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("J3").PasteSpecial(-4122)

# 4. Old "Search Page" row (now row 4): Execution Flag flips from YES to NO.
$ws.Range("C4").Value = "NO"

# 5. Old "Product Listing" row (now row 5): add the repeated
#    productCatalogPage/productDetailPage/productCatalogPage/productDetailPage/
#    productCatalogPage/productDetailPage/cartCheck columns (F:J), matching the
#    pattern already present in D5:E5.
$ws.Range("J3").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("D5:E5").Copy()
$ws.Range("F5:G5").PasteSpecial(-4122)
$ws.Range("H5:I5").PasteSpecial(-4122)
$ws.Range("F5").Value = "productCatalogPage"
$ws.Range("G5").Value = "productDetailPage"
$ws.Range("H5").Value = "productCatalogPage"
$ws.Range("I5").Value = "productDetailPage"
$ws.Range("J5").Value = "cartCheck"

# 6. Old "COD Order" row (now row 6): Execution Flag flips from YES to NO.
$ws.Range("C6").Value = "NO"

# 7. Selection as in the final workbook.
$ws.Range("A3").Select() | Out-Null
